$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing header cell text
$ws.Range("B1").Value = "Name"

# Copy the header style/formatting from B1 into the new header cells
$ws.Range("B1").Copy()
$ws.Range("C1:E1").PasteSpecial(-4122)  # xlPasteFormats

# Set the values for the new header cells
$ws.Range("C1").Value = "Potential [V]"
$ws.Range("D1").Value = "CDL [F]"
$ws.Range("E1").Value = "b [F/mV/s]"
